# users.xlsx template update:
#  - sample/help row (row 2) values refreshed to match the new
#    "refresh users / sync users" semantics:
#      A2  "alit"                               -> "ali"
#      G2  "1598656906150, 1598656906151"        -> "Staff, Students"
#      I2  "active"                              -> "enable"
#      L2  "1400/10/20 13:13:13.259"             -> "1400-10-20 13:13:13.259"
#  - active cell / selection moved from D10 to L5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ali"
$ws.Range("G2").Value = "Staff, Students"
$ws.Range("I2").Value = "enable"
$ws.Range("L2").Value = "1400-10-20 13:13:13.259"

$ws.Range("L5").Select()
